# Data_Negatif.xlsx - "fix update perhitungan btn indonesia"
#
# Update the three remaining data rows (penghasilan/pengeluaran/jangkaWaktu)
# with corrected "penghasilan" (income) figures, and remove the now-unused
# trailing example rows (5-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "penghasilan" values for the first three data rows.
$ws.Range("A2").Value = 2300000
$ws.Range("A3").Value = 3500000
$ws.Range("A4").Value = 3500000

# Remove the trailing rows that are no longer part of the sample data.
$ws.Rows("5:8").Delete() | Out-Null

# Restore the cell selection left by the editor after the edit.
$ws.Range("B13").Select() | Out-Null
